$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44519
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 3700
$ws.Range("O2").Value = 3800
$ws.Range("P2").Value = 3750
$ws.Range("Q2").Value = "$/kilo"
$ws.Range("S2").Value = 3750
$ws.Range("T2").Value = 1

# Row 3
$ws.Range("D3").Value = 44516
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 3700
$ws.Range("O3").Value = 3800
$ws.Range("P3").Value = 3750
$ws.Range("Q3").Value = "$/kilo"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 3750
$ws.Range("T3").Value = 1

# Row 4
$ws.Range("D4").Value = 44169
$ws.Range("N4").Value = 5500
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 5750
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 3833

# Row 5
$ws.Range("D5").Value = 44537
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5500
$ws.Range("P5").Value = 5250
$ws.Range("Q5").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("S5").Value = 3500
$ws.Range("T5").Value = 1.5

# Row 6
$ws.Range("D6").Value = 44159
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 6500
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6750
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 4500

# Row 7
$ws.Range("D7").Value = 44523
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 3700
$ws.Range("O7").Value = 3800
$ws.Range("P7").Value = 3750
$ws.Range("Q7").Value = "$/kilo"
$ws.Range("R7").Value = "Región del Maule"
$ws.Range("S7").Value = 3750
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44176
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 5000
$ws.Range("P8").Value = 5500
$ws.Range("S8").Value = 3667

# Row 9
$ws.Range("D9").Value = 44551

# Row 10
$ws.Range("D10").Value = 44547
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 5000
$ws.Range("O10").Value = 5500
$ws.Range("P10").Value = 5250
$ws.Range("Q10").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("S10").Value = 3500
$ws.Range("T10").Value = 1.5

# Row 11
$ws.Range("D11").Value = 44166
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 6500
$ws.Range("P11").Value = 6250
$ws.Range("Q11").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value = 4167
$ws.Range("T11").Value = 1.5

# Row 13
$ws.Range("D13").Value = 44544
$ws.Range("M13").Value = 400
$ws.Range("N13").Value = 5000
$ws.Range("O13").Value = 5500
$ws.Range("P13").Value = 5250
$ws.Range("Q13").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("S13").Value = 3500
$ws.Range("T13").Value = 1.5

# Row 14
$ws.Range("D14").Value = 44530
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 3600
$ws.Range("O14").Value = 3700
$ws.Range("P14").Value = 3650
$ws.Range("Q14").Value = "$/kilo"
$ws.Range("S14").Value = 3650
$ws.Range("T14").Value = 1

# Row 17
$ws.Range("D17").Value = 44553
$ws.Range("M17").Value = 400
$ws.Range("O17").Value = 5500
$ws.Range("P17").Value = 5250
$ws.Range("R17").Value = "Región del Maule"
$ws.Range("S17").Value = 3500
